$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Obstruction parses" -> "Obstruction parser" (bullet list item describing
#    the Obstruction parser tool).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Obstruction parses – parses obstructions relevant data from the *.fds file.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Obstruction parser – parses obstructions relevant data from the *.fds file.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Move the "QuickMark" bookmark away from the end of the paragraph that
#    talks about the /data folder contents ...
# ---------------------------------------------------------------------------
# ... and drop it right before the "Simulation has not progressed in the last
# 24h" table cell text instead. Re-adding a bookmark under the same name moves
# it, so an explicit delete of the old one is not required.
$target = $d.Content
$target.Find.Execute(
    "Simulation has not progressed in the last 24h and no stop condition is declared",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target.Collapse(1)
$d.Bookmarks.Add("QuickMark", $target) | Out-Null

# ---------------------------------------------------------------------------
# 3. Terminate two table cell sentences with a period.
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute(
    "Simulation has not progressed over a time larger than 95% of the moving standard deviation of log intervals for the last 30 logs",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r1.Collapse(0)
$r1.InsertAfter(".")

$r2 = $d.Content
$r2.Find.Execute(
    "Simulation has finished correctly",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Collapse(0)
$r2.InsertAfter(".")
